# members.xlsx edit
#  - update two "Chuyên ngành"/"Địa chỉ" values on row 2
#  - assign row 17's "Đơn vị" (dept) to "Ban Đời Sống"
#  - append a brand-new member as row 27
#  - update the saved view/selection (matches the author re-opening/re-saving
#    the sheet scrolled to the top with H6 selected)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- content edits on existing rows ---------------------------------------
$ws.Range("H2").Value = "Điện tử Viễn Thông 123"
$ws.Range("L2").Value = "Hải Dương 123"

$ws.Range("M17").Value = "Ban Đời Sống"

# --- new row 27: another member record -------------------------------------
$ws.Cells.Item(27, 1).Value = 27

# LabID ("999") must stay text, like the other alphanumeric LabIDs already in
# the sheet (B7, B17, B18, ...) - force text format, write, then drop back to
# the default "Normal" style so no stray number-format survives the save.
$ws.Cells.Item(27, 2).NumberFormat = "@"
$ws.Cells.Item(27, 2).Value = "999"
$ws.Cells.Item(27, 2).Style = "Normal"

$ws.Cells.Item(27, 3).Value = "8.png"
$ws.Cells.Item(27, 4).Value = "Ngô Xuân Hinh"
$ws.Cells.Item(27, 5).Value = "Nam"

# Ngày sinh is stored as free text everywhere in this sheet (e.g. F2:F26) -
# force text so "2022-09-01" doesn't get reinterpreted as a date serial.
$ws.Cells.Item(27, 6).NumberFormat = "@"
$ws.Cells.Item(27, 6).Value = "2022-09-01"
$ws.Cells.Item(27, 6).Style = "Normal"

$ws.Cells.Item(27, 7).Value = "1,2"
$ws.Cells.Item(27, 8).Value = "N/A"
$ws.Cells.Item(27, 9).Value = "N/A"
$ws.Cells.Item(27, 10).Value = "N/A"
$ws.Cells.Item(27, 11).Value = "N/A"
$ws.Cells.Item(27, 12).Value = "N/A"
$ws.Cells.Item(27, 13).Value = "Chưa có"
$ws.Cells.Item(27, 14).Value = "Chưa có"
$ws.Cells.Item(27, 15).Value = $false
$ws.Cells.Item(27, 16).Value = $false

# --- saved view state --------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H6").Select() | Out-Null
